# Auto-generated edit script applying numeric corrections to the Leve profit tables.
# Values were refreshed by the scheduled market-data runner (see commit message).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
# Row 33: H33,I33,J33,K33,L33,M33,N33
$ws.Range("H33").Value = 509.10526
$ws.Range("I33").Value = 297.2857
$ws.Range("J33").Value = 1102.2
$ws.Range("K33").Value = 297.2857
$ws.Range("L33").Value = 1102.2
$ws.Range("M33").Value = -68.28570000000002
$ws.Range("N33").Value = -1560.2
# Row 53: H53,I53,J53,K53,L53,M53,N53
$ws.Range("H53").Value = 1130
$ws.Range("I53").Value = 445.2857
$ws.Range("J53").Value = 1729.125
$ws.Range("K53").Value = 445.2857
$ws.Range("L53").Value = 1729.125
$ws.Range("M53").Value = 191.7143
$ws.Range("N53").Value = -3003.125
# Row 70: H70,I70,J70,K70,L70,M70,N70
$ws.Range("H70").Value = 31252062
$ws.Range("I70").Value = 2166.6667
$ws.Range("J70").Value = 50002000
$ws.Range("K70").Value = 6500.000100000001
$ws.Range("L70").Value = 150006000
$ws.Range("M70").Value = -6230.000100000001
$ws.Range("N70").Value = -150006540
# Row 73: H73,I73,J73,K73,L73,M73,N73
$ws.Range("H73").Value = 31252062
$ws.Range("I73").Value = 2166.6667
$ws.Range("J73").Value = 50002000
$ws.Range("K73").Value = 6500.000100000001
$ws.Range("L73").Value = 150006000
$ws.Range("M73").Value = -5564.000100000001
$ws.Range("N73").Value = -150007872
# Row 101: H101,J101,L101,N101
$ws.Range("H101").Value = 804.1429000000001
$ws.Range("J101").Value = 1166.3334
$ws.Range("L101").Value = 3499.0002
$ws.Range("N101").Value = -6743.0002
# Row 138: H138,I138,J138,K138,L138,M138,N138
$ws.Range("H138").Value = 27779592
$ws.Range("I138").Value = 1173.48
$ws.Range("J138").Value = 90912370
$ws.Range("K138").Value = 3520.44
$ws.Range("L138").Value = 272737110
$ws.Range("M138").Value = 1619.56
$ws.Range("N138").Value = -272747390

$ws = $wb.Worksheets("ARM")
# Row 32: H32,I32,K32,M32
$ws.Range("H32").Value = 22796600
$ws.Range("I32").Value = 23785878
$ws.Range("K32").Value = 23785878
$ws.Range("M32").Value = -23785591
# Row 34: H34,I34,K34,M34
$ws.Range("H34").Value = 212666.67
$ws.Range("I34").Value = 69000
$ws.Range("K34").Value = 69000
$ws.Range("M34").Value = -68729
# Row 63: H63,I63,K63,M63
$ws.Range("H63").Value = 3679.182
$ws.Range("I63").Value = 2299.3333
$ws.Range("K63").Value = 2299.3333
$ws.Range("M63").Value = -1613.3333
# Row 66: H66,I66,K66,M66
$ws.Range("H66").Value = 3679.182
$ws.Range("I66").Value = 2299.3333
$ws.Range("K66").Value = 11496.6665
$ws.Range("M66").Value = -8064.666499999999

$ws = $wb.Worksheets("BSM")
# Row 86: H86,I86,K86,M86
$ws.Range("H86").Value = 1557.8
$ws.Range("I86").Value = 897.5
$ws.Range("K86").Value = 897.5
$ws.Range("M86").Value = 225.5
# Row 89: H89,I89,K89,M89
$ws.Range("H89").Value = 1557.8
$ws.Range("I89").Value = 897.5
$ws.Range("K89").Value = 4487.5
$ws.Range("M89").Value = 1128.5
# Row 94: H94,I94,K94,M94
$ws.Range("H94").Value = 868
$ws.Range("I94").Value = 813.38464
$ws.Range("K94").Value = 813.38464
$ws.Range("M94").Value = -362.38464
# Row 107: H107,I107,J107,K107,L107,M107,N107
$ws.Range("H107").Value = 2247.739
$ws.Range("I107").Value = 2028.8667
$ws.Range("J107").Value = 2658.125
$ws.Range("K107").Value = 2028.8667
$ws.Range("L107").Value = 2658.125
$ws.Range("M107").Value = -108.8667
$ws.Range("N107").Value = -6498.125
# Row 134: H134,I134,K134,M134
$ws.Range("H134").Value = 7411340
$ws.Range("I134").Value = 7411340
$ws.Range("K134").Value = 22234020
$ws.Range("M134").Value = -22231485

$ws = $wb.Worksheets("CRP")
# Row 31: H31,I31,J31,K31,L31,M31,N31
$ws.Range("H31").Value = 5443.067
$ws.Range("I31").Value = 2572.1875
$ws.Range("J31").Value = 7027
$ws.Range("K31").Value = 2572.1875
$ws.Range("L31").Value = 7027
$ws.Range("M31").Value = -2277.1875
$ws.Range("N31").Value = -7617
# Row 34: H34,I34,J34,K34,L34,M34,N34
$ws.Range("H34").Value = 5443.067
$ws.Range("I34").Value = 2572.1875
$ws.Range("J34").Value = 7027
$ws.Range("K34").Value = 2572.1875
$ws.Range("L34").Value = 7027
$ws.Range("M34").Value = -2370.1875
$ws.Range("N34").Value = -7431

$ws = $wb.Worksheets("CUL")
# Row 14: H14,I14,K14,M14
$ws.Range("H14").Value = 18503.834
$ws.Range("I14").Value = 18503.834
$ws.Range("K14").Value = 55511.50199999999
$ws.Range("M14").Value = -55338.50199999999
# Row 39: H39,J39,L39,N39
$ws.Range("H39").Value = 4210.3
$ws.Range("J39").Value = 4289.7144
$ws.Range("L39").Value = 12869.1432
$ws.Range("N39").Value = -13457.1432
# Row 113: H113,J113,L113,N113
$ws.Range("H113").Value = 3714.1428
$ws.Range("J113").Value = 4166.5
$ws.Range("L113").Value = 12499.5
$ws.Range("N113").Value = -16839.5
# Row 134: H134,I134,K134,M134
$ws.Range("H134").Value = 2946.05
$ws.Range("I134").Value = 2946.05
$ws.Range("K134").Value = 8838.150000000001
$ws.Range("M134").Value = -3768.150000000001

$ws = $wb.Worksheets("GSM")
# Row 20: H20,I20,J20,K20,L20,M20,N20
$ws.Range("H20").Value = 18599.4
$ws.Range("I20").Value = 3000
$ws.Range("J20").Value = 22499.25
$ws.Range("K20").Value = 3000
$ws.Range("L20").Value = 22499.25
$ws.Range("M20").Value = -2755
$ws.Range("N20").Value = -22989.25
# Row 24: H24,J24,L24,N24
$ws.Range("H24").Value = 4032600
$ws.Range("J24").Value = 49998
$ws.Range("L24").Value = 49998
$ws.Range("N24").Value = -50344
# Row 33: H33,J33,L33,N33
$ws.Range("H33").Value = 22500000
$ws.Range("J33").Value = 22500000
$ws.Range("L33").Value = 22500000
$ws.Range("N33").Value = -22500504
# Row 107: H107,I107,K107,M107
$ws.Range("H107").Value = 980
$ws.Range("I107").Value = 841.0909
$ws.Range("K107").Value = 841.0909
$ws.Range("M107").Value = 1078.9091
# Row 122: H122,I122,J122,K122,L122,M122,N122
$ws.Range("H122").Value = 4529.4287
$ws.Range("I122").Value = 4224.75
$ws.Range("J122").Value = 4935.6665
$ws.Range("K122").Value = 12674.25
$ws.Range("L122").Value = 14806.9995
$ws.Range("M122").Value = -10224.25
$ws.Range("N122").Value = -19706.9995
# Row 126: H126,I126,K126,M126
$ws.Range("H126").Value = 2853.75
$ws.Range("I126").Value = 2596.5
$ws.Range("K126").Value = 7789.5
$ws.Range("M126").Value = -5319.5
# Row 132: H132,I132,J132,K132,L132,M132,N132
$ws.Range("H132").Value = 3412.2083
$ws.Range("I132").Value = 3386.652
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 10159.956
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -7629.956
$ws.Range("N132").Value = -17060

$ws = $wb.Worksheets("LTW")
# Row 82: H82,I82,J82,K82,L82,M82,N82
$ws.Range("H82").Value = 2559.6365
$ws.Range("I82").Value = 2018.6154
$ws.Range("J82").Value = 3341.111
$ws.Range("K82").Value = 2018.6154
$ws.Range("L82").Value = 3341.111
$ws.Range("M82").Value = -1657.6154
$ws.Range("N82").Value = -4063.111
# Row 85: H85,I85,J85,K85,L85,M85,N85
$ws.Range("H85").Value = 2559.6365
$ws.Range("I85").Value = 2018.6154
$ws.Range("J85").Value = 3341.111
$ws.Range("K85").Value = 2018.6154
$ws.Range("L85").Value = 3341.111
$ws.Range("M85").Value = -770.6153999999999
$ws.Range("N85").Value = -5837.111
# Row 100: H100,I100,K100,M100
$ws.Range("H100").Value = 1300.8182
$ws.Range("I100").Value = 1145.4445
$ws.Range("K100").Value = 1145.4445
$ws.Range("M100").Value = -604.4445000000001
# Row 132: H132,I132,K132,M132
$ws.Range("H132").Value = 4926.2
$ws.Range("I132").Value = 4776.385
$ws.Range("K132").Value = 14329.155
$ws.Range("M132").Value = -11799.155

$ws = $wb.Worksheets("WVR")
# Row 31: H31,J31,L31,N31
$ws.Range("H31").Value = 34685
$ws.Range("J31").Value = 42019
$ws.Range("L31").Value = 42019
$ws.Range("N31").Value = -42715
# Row 41: H41,J41,L41,N41
$ws.Range("H41").Value = 7999.5
$ws.Range("J41").Value = 7999.5
$ws.Range("L41").Value = 7999.5
$ws.Range("N41").Value = -8779.5
# Row 62: H62,J62,L62,N62
$ws.Range("H62").Value = 5062.75
$ws.Range("J62").Value = 5667
$ws.Range("L62").Value = 5667
$ws.Range("N62").Value = -6915
# Row 65: H65,J65,L65,N65
$ws.Range("H65").Value = 5062.75
$ws.Range("J65").Value = 5667
$ws.Range("L65").Value = 28335
$ws.Range("N65").Value = -34575
# Row 122: H122,I122,J122,K122,L122,M122,N122
$ws.Range("H122").Value = 4912.231
$ws.Range("I122").Value = 4203.476
$ws.Range("J122").Value = 7889
$ws.Range("K122").Value = 12610.428
$ws.Range("L122").Value = 23667
$ws.Range("M122").Value = -10160.428
$ws.Range("N122").Value = -28567
# Row 136: H136,I136,K136,M136
$ws.Range("H136").Value = 1350.0646
$ws.Range("I136").Value = 1253.7222
$ws.Range("K136").Value = 3761.1666
$ws.Range("M136").Value = -1211.1666
